# semana 46 de 2025
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 1

# Row 5
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 0.12

# Row 6
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 0.1

# Row 7
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 0

# Row 9
$ws.Range("C9").Value = 41
$ws.Range("D9").Value = 33

# Row 10
$ws.Range("C10").Value = 1
$ws.Range("E10").Value = 0.37

# Row 11
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0.37

# Row 12
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 9
$ws.Range("E12").Value = 0

# Row 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.37

# Row 14
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0.37

# Row 16
$ws.Range("C16").Value = 0
$ws.Range("E16").Value = 1

# Row 17
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 0.09

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0.37

# Row 19
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 0.09

# Row 21
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0

# Row 25
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0.03

# Row 26
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0.37

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0

# Row 29
$ws.Range("C29").Value = 0
$ws.Range("E29").Value = 1

# Row 31
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 0.18

# Row 33
$ws.Range("C33").Value = 7
$ws.Range("D33").Value = 4
$ws.Range("E33").Value = 0.09

# Row 34
$ws.Range("C34").Value = 11
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0

# Row 35
$ws.Range("C35").Value = 6
$ws.Range("D35").Value = 10
$ws.Range("E35").Value = 0.04
